$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.729.53"
$ws.Range("E2").Value = "  -3.25%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.098.48"
$ws.Range("E3").Value = "  -2.31%  "
# Row 4
$ws.Range("E4").Value = "  -0.29%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.27"
$ws.Range("E5").Value = "  -2.12%  "
# Row 6
$ws.Range("E6").Value = "  -0.15%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5150"
$ws.Range("E7").Value = "  -2.55%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4414"
$ws.Range("E8").Value = "  -3.48%  "
# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09297"
$ws.Range("E9").Value = "  +1.15%  "
# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.65"
$ws.Range("E10").Value = "  -3.06%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.173"
$ws.Range("E11").Value = "  -1.05%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.90"
$ws.Range("E12").Value = "  -0.26%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.104.46"
$ws.Range("E13").Value = "  -1.99%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.291"
$ws.Range("E14").Value = "  +1.43%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.763"
$ws.Range("E15").Value = "  -2.20%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.74"
$ws.Range("E16").Value = "  -2.72%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001152"
$ws.Range("E17").Value = "  -2.43%  "
# Row 18
$ws.Range("E18").Value = "  -0.23%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.90"
$ws.Range("E19").Value = "  +6.52%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06644"
$ws.Range("E20").Value = "  -1.30%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -0.13%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.202"
$ws.Range("E22").Value = "  -2.78%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.787.06"
$ws.Range("E23").Value = "  -3.34%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.61"
$ws.Range("E24").Value = "  -2.19%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  -3.02%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.358.33"
$ws.Range("E26").Value = "  -0.80%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.94"
$ws.Range("E27").Value = "  -2.92%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.529"
$ws.Range("E28").Value = "  -3.92%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.14"
$ws.Range("E29").Value = "  -1.88%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.20"
$ws.Range("E30").Value = "  -2.74%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.133"
$ws.Range("E31").Value = "  -7.29%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("E32").Value = "  -3.04%  "
# Row 33
$ws.Range("E33").Value = "  -1.07%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.185"
$ws.Range("E34").Value = "  -3.72%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.943"
$ws.Range("E35").Value = "  -1.88%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.115"
$ws.Range("E36").Value = "  -0.49%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.38"
$ws.Range("E37").Value = "  -0.72%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("E38").Value = "  -3.37%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06730"
$ws.Range("E39").Value = "  -2.92%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.48"
$ws.Range("E40").Value = "  -1.95%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6862"
$ws.Range("E41").Value = "  -1.99%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2231"
$ws.Range("E42").Value = "  -4.57%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.298"
$ws.Range("E43").Value = "  +1.86%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6640"
$ws.Range("E44").Value = "  +2.29%  "
# Row 45
$ws.Range("E45").Value = "  -3.61%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.326"
$ws.Range("E46").Value = "  -1.79%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.621"
$ws.Range("E47").Value = "  -3.66%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000347"
$ws.Range("E48").Value = "  -6.31%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.222"
$ws.Range("E49").Value = "  -3.10%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.41"
$ws.Range("E50").Value = "  -1.24%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3347"
$ws.Range("E51").Value = "  +0.13%  "
